# Natmi following Dr Hou advice
# Update Wnt1-Fzd4 LR-pair sheet: recompute rows 2-4 (now sent from "ECs") and
# append new rows 5-7 for the "FAPs" sending cluster with updated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs", "Wnt1", "Fzd4", "ECs", 1, "0.3333333333333333", "0.01948966666666667", "0.058469", "0.07096062449330311", "0.07096062449330311", 3, 1, "29.22839533333334", "87.685186", "0.4452121315669252", "0.4988525987881078", "0.5696516822482223", "5.126865140234", "0.03159253088798364", "0.03539889194011131")
    ,@("ECs", "Wnt1", "Fzd4", "FAPs", 1, "0.3333333333333333", "0.01948966666666667", "0.058469", "0.07096062449330311", "0.07096062449330311", 3, 1, "15.24435933333334", "45.73307800000001", "0.232204800700274", "0.2601815181287206", "0.2971074819535556", "2.673967337582", "0.01647739766803443", "0.01846264300802968")
    ,@("ECs", "Wnt1", "Fzd4", "sCs", 1, "0.3333333333333333", "0.01948966666666667", "0.058469", "0.07096062449330311", "0.07096062449330311", 2, 1, "21.177737", "42.355474", "0.3225830677328006", "0.2409658830831714", "0.4127470348843333", "2.476482209306", "0.02289069593728503", "0.0170990895451621")
    ,@("FAPs", "Wnt1", "Fzd4", "ECs", 3, 1, "0.255165", "0.7654949999999999", "0.9290393755066968", "0.9290393755066969", 3, 1, "29.22839533333334", "87.685186", "0.4452121315669252", "0.4988525987881078", "7.458063495229999", "67.12257145706999", "0.4136196006789415", "0.4634537068479965")
    ,@("FAPs", "Wnt1", "Fzd4", "FAPs", 3, 1, "0.255165", "0.7654949999999999", "0.9290393755066968", "0.9290393755066969", 3, 1, "15.24435933333334", "45.73307800000001", "0.232204800700274", "0.2601815181287206", "3.88982694929", "35.00844254361", "0.2157274030322396", "0.241718875120691")
    ,@("FAPs", "Wnt1", "Fzd4", "sCs", 3, 1, "0.255165", "0.7654949999999999", "0.9290393755066968", "0.9290393755066969", 2, 1, "21.177737", "42.355474", "0.3225830677328006", "0.2409658830831714", "5.403817261605", "32.42290356962999", "0.2996923717955156", "0.2238667935380093")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}

Write-Output "Updated Sheet1 rows 2-7 (added 3 new rows for FAPs sending cluster)"
